# Apply the DNB 2025-05 "increase merchant diversity" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet tab.
$ws.Name = "Sheet"

# 2) Widen the date number format to include a time component (applies
#    to the whole date column so every row ends up on the same style).
$ws.Range("A2:A20").NumberFormat = "yyyy-mm-dd h:mm:ss"

# 3) Replace the transaction rows (2-14) with the new data set, and
#    append the newly-added rows (15-20).

# Row 2
$ws.Range("A2").Value = 45808
$ws.Range("B2").Value = "SPOTIFY"
$ws.Range("F2").Value = 129

# Row 3
$ws.Range("A3").Value = 45807
$ws.Range("B3").Value = "KIWI GRØNLAND"
$ws.Range("F3").Value = 378.9

# Row 4
$ws.Range("A4").Value = 45805
$ws.Range("B4").Value = "NETFLIX.COM"
$ws.Range("F4").Value = 179

# Row 5
$ws.Range("A5").Value = 45804
$ws.Range("B5").Value = "BLI VAKKER MAJORSTUEN"
$ws.Range("F5").Value = 456

# Row 6
$ws.Range("A6").Value = 45802
$ws.Range("B6").Value = "REMA1000 MAJORSTUEN"
$ws.Range("F6").Value = 567.3

# Row 7
$ws.Range("A7").Value = 45801
$ws.Range("B7").Value = "PIZZABAKEREN TORSHOV"
$ws.Range("F7").Value = 349

# Row 8
$ws.Range("A8").Value = 45799
$ws.Range("B8").Value = "VINMONOPOLET STORO"
$ws.Range("F8").Value = 678

# Row 9
$ws.Range("A9").Value = 45797
$ws.Range("B9").Value = "STARBUCKS MAJORSTUEN"
$ws.Range("F9").Value = 85

# Row 10
$ws.Range("A10").Value = 45795
$ws.Range("B10").Value = "BYGGMAX ALNA"
$ws.Range("F10").Value = 2345

# Row 11
$ws.Range("A11").Value = 45794
$ws.Range("B11").Value = "SAS PLUS ARLANDA"
$ws.Range("F11").Value = 3499

# Row 12
$ws.Range("A12").Value = 45792
$ws.Range("B12").Value = "NORWEGIAN GARDERMOEN"
$ws.Range("F12").Value = 1899

# Row 13 (was "Innbetaling"/E13 - now a regular expense row; clear E13)
$ws.Range("A13").Value = 45791
$ws.Range("B13").Value = "PLANTASJEN SINSEN"
$ws.Range("E13").ClearContents()
$ws.Range("F13").Value = 789

# Row 14
$ws.Range("A14").Value = 45789
$ws.Range("B14").Value = "MENY GRØNLAND"
$ws.Range("F14").Value = 489

# Row 15 (new)
$ws.Range("A15").Value = 45787
$ws.Range("B15").Value = "CLAES OHLSON STORO"
$ws.Range("F15").Value = 234

# Row 16 (new)
$ws.Range("A16").Value = 45785
$ws.Range("B16").Value = "EUROPRIS TORSHOV"
$ws.Range("F16").Value = 189

# Row 17 (new)
$ws.Range("A17").Value = 45783
$ws.Range("B17").Value = "KAFFEBRENNERIET MAJORSTUEN"
$ws.Range("F17").Value = 89

# Row 18 (new - the relocated "Innbetaling" income row)
$ws.Range("A18").Value = 45782
$ws.Range("B18").Value = "Innbetaling"
$ws.Range("E18").Value = 15000

# Row 19 (new)
$ws.Range("A19").Value = 45780
$ws.Range("B19").Value = "COOP EXTRA SAGENE"
$ws.Range("F19").Value = 678.5

# Row 20 (new)
$ws.Range("A20").Value = 45778
$ws.Range("B20").Value = "TEKNIKMAGASINET STORO"
$ws.Range("F20").Value = 299
